$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'travel fitness'
$ws.Range("A2").Value = 'travel leg support'
$ws.Range("A3").Value = 'travel leggings for women'
$ws.Range("A4").Value = 'travel package for women'
$ws.Range("A5").Value = 'travel tights'
$ws.Range("A6").Value = 'travel waist'
$ws.Range("A7").Value = 'travel workout'
$ws.Range("A8").Value = 'tree leggings for women'
$ws.Range("A9").Value = 'trick bikes'
$ws.Range("A10").Value = 'trick track game'
$ws.Range("A11").Value = 'tricky nearly god'
$ws.Range("A12").Value = 'trouble in a tight dress'
$ws.Range("A13").Value = 'trousers with braces'
$ws.Range("A14").Value = 'true flight size'
$ws.Range("A15").Value = 'true joint'
$ws.Range("A16").Value = 'true recovery'
$ws.Range("A17").Value = 'true tear device'
$ws.Range("A18").Value = 'true tech pants'
$ws.Range("A19").Value = 'trx bike'
$ws.Range("A20").Value = 'trx on the go'
$ws.Range("A21").Value = 'trx stand'
$ws.Range("A22").Value = 'trx to go'
$ws.Range("A23").Value = 'trx total recovery'
$ws.Range("A24").Value = 'tsla yoga'
$ws.Range("A25").Value = 'tummy compression band'
$ws.Range("A26").Value = 'tummy compression panties'
$ws.Range("A27").Value = 'tummy compression tank'
$ws.Range("A28").Value = 'tummy control undies'
$ws.Range("A29").Value = 'tummy panel leggings'
$ws.Range("A30").Value = 'tummy support leggings'
$ws.Range("A31").Value = 'tummy support panties'
$ws.Range("A32").Value = 'tummy tape'
$ws.Range("A33").Value = 'tummy tight'
$ws.Range("A34").Value = 'tummy tight for her'
$ws.Range("A35").Value = 'tummy warmer for women'
$ws.Range("A36").Value = 'two toned pants'
$ws.Range("A37").Value = 'two toned workout pants'
$ws.Range("A38").Value = 'ua capris women'
$ws.Range("A39").Value = 'ultimate frisbee apparel'
$ws.Range("A40").Value = 'ultimate frisbee clothes'
$ws.Range("A41").Value = 'ultimate frisbee clothing'
$ws.Range("A42").Value = 'ultimate frisbee equipment'
$ws.Range("A43").Value = 'ultimate frisbee game'
$ws.Range("A44").Value = 'ultimate frisbee gear'
$ws.Range("A45").Value = 'ultimate frisbee pants'
$ws.Range("A46").Value = 'ultimate frisbee shorts'
$ws.Range("A47").Value = 'ultimate frisbee training'
$ws.Range("A48").Value = 'ultimate workout wear'
$ws.Range("A49").Value = 'ultimate yoga leggings'
$ws.Range("A50").Value = 'ultra black leggings'
$ws.Range("A51").Value = 'ultra compression leggings'
$ws.Range("A52").Value = 'ultramarathon gear'
$ws.Range("A53").Value = 'underwear for running women'
$ws.Range("A54").Value = 'underwear pocket women'
$ws.Range("A55").Value = 'underwear running women'
$ws.Range("A56").Value = 'underwear women days of the week'
$ws.Range("A57").Value = 'underwear women for periods'
$ws.Range("A58").Value = 'underwear women period'
$ws.Range("A59").Value = 'underwear women running'
$ws.Range("A60").Value = 'undies for two'
$ws.Range("A61").Value = 'unicorm leggings'
$ws.Range("A62").Value = 'up right freezers'
$ws.Range("A63").Value = 'up tight'
$ws.Range("A64").Value = 'upf leggings'
$ws.Range("A65").Value = 'upf leggings women'
$ws.Range("A66").Value = 'upper back compression'
$ws.Range("A67").Value = 'upper leg joint pain'
$ws.Range("A68").Value = 'upper thigh compression'
$ws.Range("A69").Value = 'used generator'
$ws.Range("A70").Value = 'uv joint'
$ws.Range("A71").Value = 'uv pants women'
$ws.Range("A72").Value = 'uv total recovery'
$ws.Range("A73").Value = 'uva basketball apparel'
$ws.Range("A74").Value = 'uva clothing for women'
$ws.Range("A75").Value = 'uvb light therapy'
$ws.Range("A76").Value = 'uvb light therapy for skin'
$ws.Range("A77").Value = 'ventilator machine'
$ws.Range("A78").Value = 'ventilator medical'
$ws.Range("A79").Value = 'ventilator medical machine'
$ws.Range("A80").Value = 'vertigo band'
$ws.Range("A81").Value = 'vertigo clothing women'
$ws.Range("A82").Value = 'vertigo help'
$ws.Range("A83").Value = 'vertigo inhaler'
$ws.Range("A84").Value = 'vertigo womens clothing'
$ws.Range("A85").Value = 'vesture hot pack replacement'
$ws.Range("A86").Value = 'victoria leggings'
$ws.Range("A87").Value = 'victoria secret black yoga pants'
$ws.Range("A88").Value = 'victoria secret capri leggings'
$ws.Range("A89").Value = 'victoria secret gift baskets for women'
$ws.Range("A90").Value = 'victoria secret high waist'
$ws.Range("A91").Value = 'victoria secret leggings'
$ws.Range("A92").Value = 'victoria secret leggings cheap'
$ws.Range("A93").Value = 'victoria secret leggings for women'
$ws.Range("A94").Value = 'victoria secret leggings with pockets'
$ws.Range("A95").Value = 'victoria secret leggings xs'
$ws.Range("A96").Value = 'victoria secret pink leggins'
$ws.Range("A97").Value = 'victoria secret rainbow'
$ws.Range("A98").Value = 'victoria secret running shorts'
$ws.Range("A99").Value = 'victoria secret shorts'
$ws.Range("A100").Value = 'victoria secret sport pants'
